# Apply updated TPM-derived values (new ligand/receptor expression computations)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("G2").Value = 90.25004833333332
$ws.Range("H2").Value = 270.750145
$ws.Range("I2").Value = 0.8928575650827933
$ws.Range("J2").Value = 0.8928575650827932
$ws.Range("M2").Value = 32.736679
$ws.Range("N2").Value = 98.210037
$ws.Range("O2").Value = 0.8346853755332739
$ws.Range("P2").Value = 0.834685375533274
$ws.Range("Q2").Value = 2954.486862022818
$ws.Range("R2").Value = 26590.38175820536
$ws.Range("S2").Value = 0.7452551520088559
$ws.Range("T2").Value = 0.7452551520088559

# Row 3
$ws.Range("G3").Value = 90.25004833333332
$ws.Range("H3").Value = 270.750145
$ws.Range("I3").Value = 0.8928575650827933
$ws.Range("J3").Value = 0.8928575650827932
$ws.Range("O3").Value = 0.001766029048926899
$ws.Range("P3").Value = 0.0017660290489269
$ws.Range("Q3").Value = 6.251109431109444
$ws.Range("R3").Value = 56.259984879985
$ws.Range("S3").Value = 0.001576812396490352
$ws.Range("T3").Value = 0.001576812396490353

# Row 4
$ws.Range("G4").Value = 90.25004833333332
$ws.Range("H4").Value = 270.750145
$ws.Range("I4").Value = 0.8928575650827933
$ws.Range("J4").Value = 0.8928575650827932
$ws.Range("M4").Value = 0.5119106666666666
$ws.Range("N4").Value = 1.535732
$ws.Range("O4").Value = 0.01305215923234471
$ws.Range("P4").Value = 0.01305215923234471
$ws.Range("Q4").Value = 46.19996240901555
$ws.Range("R4").Value = 415.7996616811399
$ws.Range("S4").Value = 0.0116537191112642
$ws.Range("T4").Value = 0.0116537191112642

# Row 5
$ws.Range("G5").Value = 90.25004833333332
$ws.Range("H5").Value = 270.750145
$ws.Range("I5").Value = 0.8928575650827933
$ws.Range("J5").Value = 0.8928575650827932
$ws.Range("M5").Value = 5.902527666666667
$ws.Range("N5").Value = 17.707583
$ws.Range("O5").Value = 0.1504964361854544
$ws.Range("P5").Value = 0.1504964361854544
$ws.Range("Q5").Value = 532.7034072055038
$ws.Range("R5").Value = 4794.330664849535
$ws.Range("S5").Value = 0.1343718815661828
$ws.Range("T5").Value = 0.1343718815661828

# Row 6
$ws.Range("I6").Value = 0.05133510428912089
$ws.Range("J6").Value = 0.05133510428912089
$ws.Range("M6").Value = 32.736679
$ws.Range("N6").Value = 98.210037
$ws.Range("O6").Value = 0.8346853755332739
$ws.Range("P6").Value = 0.834685375533274
$ws.Range("Q6").Value = 169.8690777948607
$ws.Range("R6").Value = 1528.821700153746
$ws.Range("S6").Value = 0.04284866080160465
$ws.Range("T6").Value = 0.04284866080160465

# Row 7
$ws.Range("I7").Value = 0.05133510428912089
$ws.Range("J7").Value = 0.05133510428912089
$ws.Range("O7").Value = 0.001766029048926899
$ws.Range("P7").Value = 0.0017660290489269
$ws.Range("S7").Value = 0.00009065928540427934
$ws.Range("T7").Value = 0.00009065928540427937

# Row 8
$ws.Range("I8").Value = 0.05133510428912089
$ws.Range("J8").Value = 0.05133510428912089
$ws.Range("M8").Value = 0.5119106666666666
$ws.Range("N8").Value = 1.535732
$ws.Range("O8").Value = 0.01305215923234471
$ws.Range("P8").Value = 0.01305215923234471
$ws.Range("Q8").Value = 2.656280218895111
$ws.Range("R8").Value = 23.906521970056
$ws.Range("S8").Value = 0.0006700339553906278
$ws.Range("T8").Value = 0.0006700339553906278

# Row 9
$ws.Range("I9").Value = 0.05133510428912089
$ws.Range("J9").Value = 0.05133510428912089
$ws.Range("M9").Value = 5.902527666666667
$ws.Range("N9").Value = 17.707583
$ws.Range("O9").Value = 0.1504964361854544
$ws.Range("P9").Value = 0.1504964361854544
$ws.Range("Q9").Value = 30.62793667602378
$ws.Range("R9").Value = 275.651430084214
$ws.Range("S9").Value = 0.007725750246721328
$ws.Range("T9").Value = 0.00772575024672133

# Row 10
$ws.Range("G10").Value = 5.380476000000001
$ws.Range("H10").Value = 16.141428
$ws.Range("I10").Value = 0.05322987398968605
$ws.Range("J10").Value = 0.05322987398968604
$ws.Range("M10").Value = 32.736679
$ws.Range("N10").Value = 98.210037
$ws.Range("O10").Value = 0.8346853755332739
$ws.Range("P10").Value = 0.834685375533274
$ws.Range("Q10").Value = 176.138915679204
$ws.Range("R10").Value = 1585.250241112836
$ws.Range("S10").Value = 0.04443019736066995
$ws.Range("T10").Value = 0.04443019736066995

# Row 11
$ws.Range("G11").Value = 5.380476000000001
$ws.Range("H11").Value = 16.141428
$ws.Range("I11").Value = 0.05322987398968605
$ws.Range("J11").Value = 0.05322987398968604
$ws.Range("O11").Value = 0.001766029048926899
$ws.Range("P11").Value = 0.0017660290489269
$ws.Range("Q11").Value = 0.372675083156
$ws.Range("R11").Value = 3.354075748404
$ws.Range("S11").Value = 0.00009400550373650393
$ws.Range("T11").Value = 0.00009400550373650394

# Row 12
$ws.Range("G12").Value = 5.380476000000001
$ws.Range("H12").Value = 16.141428
$ws.Range("I12").Value = 0.05322987398968605
$ws.Range("J12").Value = 0.05322987398968604
$ws.Range("M12").Value = 0.5119106666666666
$ws.Range("N12").Value = 1.535732
$ws.Range("O12").Value = 0.01305215923234471
$ws.Range("P12").Value = 0.01305215923234471
$ws.Range("Q12").Value = 2.754323056144
$ws.Range("R12").Value = 24.788907505296
$ws.Range("S12").Value = 0.0006947647912310264
$ws.Range("T12").Value = 0.0006947647912310264

# Row 13
$ws.Range("G13").Value = 5.380476000000001
$ws.Range("H13").Value = 16.141428
$ws.Range("I13").Value = 0.05322987398968605
$ws.Range("J13").Value = 0.05322987398968604
$ws.Range("M13").Value = 5.902527666666667
$ws.Range("N13").Value = 17.707583
$ws.Range("O13").Value = 0.1504964361854544
$ws.Range("P13").Value = 0.1504964361854544
$ws.Range("Q13").Value = 31.758408449836
$ws.Range("R13").Value = 285.825676048524
$ws.Range("S13").Value = 0.008010906334048565
$ws.Range("T13").Value = 0.008010906334048567

# Row 14
$ws.Range("E14").Value = 3
$ws.Range("F14").Value = 1
$ws.Range("G14").Value = 0.2605293333333333
$ws.Range("H14").Value = 0.7815879999999999
$ws.Range("I14").Value = 0.002577456638399696
$ws.Range("J14").Value = 0.002577456638399696
$ws.Range("M14").Value = 32.736679
$ws.Range("N14").Value = 98.210037
$ws.Range("O14").Value = 0.8346853755332739
$ws.Range("P14").Value = 0.834685375533274
$ws.Range("Q14").Value = 8.528865155417334
$ws.Range("R14").Value = 76.75978639875599
$ws.Range("S14").Value = 0.002151365362143381
$ws.Range("T14").Value = 0.002151365362143381

# Row 15
$ws.Range("E15").Value = 3
$ws.Range("F15").Value = 1
$ws.Range("G15").Value = 0.2605293333333333
$ws.Range("H15").Value = 0.7815879999999999
$ws.Range("I15").Value = 0.002577456638399696
$ws.Range("J15").Value = 0.002577456638399696
$ws.Range("O15").Value = 0.001766029048926899
$ws.Range("P15").Value = 0.0017660290489269
$ws.Range("Q15").Value = 0.01804539058711111
$ws.Range("R15").Value = 0.162408515284
$ws.Range("S15").Value = 0.000004551863295763338
$ws.Range("T15").Value = 0.000004551863295763339

# Row 16
$ws.Range("E16").Value = 3
$ws.Range("F16").Value = 1
$ws.Range("G16").Value = 0.2605293333333333
$ws.Range("H16").Value = 0.7815879999999999
$ws.Range("I16").Value = 0.002577456638399696
$ws.Range("J16").Value = 0.002577456638399696
$ws.Range("M16").Value = 0.5119106666666666
$ws.Range("N16").Value = 1.535732
$ws.Range("O16").Value = 0.01305215923234471
$ws.Range("P16").Value = 0.01305215923234471
$ws.Range("Q16").Value = 0.1333677447128889
$ws.Range("R16").Value = 1.200309702416
$ws.Range("S16").Value = 0.00003364137445885676
$ws.Range("T16").Value = 0.00003364137445885676

# Row 17
$ws.Range("E17").Value = 3
$ws.Range("F17").Value = 1
$ws.Range("G17").Value = 0.2605293333333333
$ws.Range("H17").Value = 0.7815879999999999
$ws.Range("I17").Value = 0.002577456638399696
$ws.Range("J17").Value = 0.002577456638399696
$ws.Range("M17").Value = 5.902527666666667
$ws.Range("N17").Value = 17.707583
$ws.Range("O17").Value = 0.1504964361854544
$ws.Range("P17").Value = 0.1504964361854544
$ws.Range("Q17").Value = 1.537781597978222
$ws.Range("R17").Value = 13.840034381804
$ws.Range("S17").Value = 0.0003878980385016958
$ws.Range("T17").Value = 0.0003878980385016958
